$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.774.86'
$ws.Range("E2").Value = '  +7.08%  '

$ws.Range("D3").Value = '2.401.99'
$ws.Range("E3").Value = '  +4.30%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '115.82'
$ws.Range("E5").Value = '  +10.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '319.84'
$ws.Range("E6").Value = '  +3.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.639'
$ws.Range("E7").Value = '  +3.30%  '

$ws.Range("E8").Value = '  -0.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.632'
$ws.Range("E9").Value = '  +4.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.97'
$ws.Range("E10").Value = '  +8.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0933'
$ws.Range("E11").Value = '  +3.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.75'
$ws.Range("E12").Value = '  +5.70%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.111'
$ws.Range("E13").Value = '  +2.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.02'
$ws.Range("E14").Value = '  +3.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.06'
$ws.Range("E15").Value = '  +4.51%  '

$ws.Range("D16").Value = '2.766.40'
$ws.Range("E16").Value = '  -0.58%  '

$ws.Range("D17").Value = '2.411.44'
$ws.Range("E17").Value = '  +4.70%  '

$ws.Range("D18").Value = '45.736.74'
$ws.Range("E18").Value = '  +6.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.56'
$ws.Range("E19").Value = '  +3.29%  '

$ws.Range("E20").Value = '  +3.69%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.74'
$ws.Range("E21").Value = '  +1.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.31'
$ws.Range("E22").Value = '  +2.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.57'
$ws.Range("E23").Value = '  +4.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '265.87'
$ws.Range("E24").Value = '  -0.99%  '

$ws.Range("E25").Value = '  +7.93%  '

$ws.Range("E26").Value = '  -0.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.66'
$ws.Range("E27").Value = '  +6.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.42'
$ws.Range("E28").Value = '  +5.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.36'
$ws.Range("E29").Value = '  +2.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '40.69'
$ws.Range("E30").Value = '  +12.44%  '

$ws.Range("E31").Value = '  +16.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.90'
$ws.Range("E32").Value = '  +2.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '173.33'
$ws.Range("E33").Value = '  +5.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.96'
$ws.Range("E34").Value = '  +12.57%  '

$ws.Range("E35").Value = '  +1.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.03'
$ws.Range("E36").Value = '  +10.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.121'
$ws.Range("E37").Value = '  +8.81%  '

$ws.Range("E38").Value = '  +17.37%  '

$ws.Range("E39").Value = '  +10.79%  '

$ws.Range("E40").Value = '  +5.08%  '

$ws.Range("E41").Value = '  +11.77%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.244'
$ws.Range("E42").Value = '  +7.61%  '

$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.83'
$ws.Range("E43").Value = '  +12.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.99'
$ws.Range("E44").Value = '  -9.25%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '72.61'
$ws.Range("E45").Value = '  +2.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.69'
$ws.Range("E46").Value = '  +14.15%  '

$ws.Range("E47").Value = '  -0.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.91'
$ws.Range("E48").Value = '  +14.94%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.87'
$ws.Range("E49").Value = '  +5.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.50'
$ws.Range("E50").Value = '  +9.78%  '

$ws.Range("E51").Value = '  +12.87%  '
